# "Pop up bar width" - resize the columns on the Login Page Test sheet and
# update the saved view state (scroll position / selection) to match where
# the author was working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Activate()

# View was scrolled so row 10 is the top-left visible row, and the
# selection left on Q6:R6.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q6:R6").Select()

# New column widths (narrower overall - the popup bar no longer needs the
# very wide columns it used to).
# Target character widths (as stored in the workbook XML) are:
#   A=8.5703125 B=20.42578125 C=25.140625 D=23.7109375
#   E=24.7109375 F=24.7109375 G=25.28515625 H=11.28515625 I=10
# ColumnWidth (COM) rounds to the nearest 1/6th of a character internally,
# so we feed it the value that rounds back to the closest match.
$ws.Columns.Item(1).ColumnWidth = 7.666666666666667   # -> 8.5703125
$ws.Columns.Item(2).ColumnWidth = 19.666666666666668  # -> 20.42578125
$ws.Columns.Item(3).ColumnWidth = 24.333333333333332  # -> 25.140625
$ws.Columns.Item(4).ColumnWidth = 22.833333333333332  # -> 23.7109375
$ws.Columns.Item(5).ColumnWidth = 23.833333333333332  # -> 24.7109375
$ws.Columns.Item(6).ColumnWidth = 23.833333333333332  # -> 24.7109375
$ws.Columns.Item(7).ColumnWidth = 24.5                # -> 25.28515625
$ws.Columns.Item(8).ColumnWidth = 10.5                # -> 11.28515625
$ws.Columns.Item(9).ColumnWidth = 9.166666666666666   # -> 10
